$d = $word.ActiveDocument

$d.Content.Find.Execute("Is altruist", $true, $false, $false, $false, $false, $true, 1, $false, "Is willing to donate to reforestation project", 2)

$d.Content.Find.Execute('wtp == "Yes" / donation > median(donation)', $true, $false, $false, $false, $false, $true, 1, $false, "donation > median(donation)", 2)

$d.Content.Find.Execute("Is willing to adapt to climate change", $true, $false, $false, $false, $false, $true, 1, $false, "Is willing to adopt climate friendly behavior", 2)
